$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") — reuse the same style as the
# existing header cells (e.g. H1: bold, centered, bordered) by copying it.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells for rows 2-3
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9
